$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("median")

$data = @(
    @(23.54, 37.18, 161.48),
    @(34.32, 25.52, 12.76),
    @(27.94, 43.34, 33.66),
    @(26.18, 36.52, 35.86),
    @(27.5,  16.06, 106.48),
    @($null, 57.97, 29.48),
    @($null, 32.56, 97.24),
    @($null, $null, 19.14),
    @($null, $null, 45.54),
    @($null, $null, 21.34),
    @($null, $null, 16.94),
    @($null, $null, 34.32),
    @($null, $null, 164.12),
    @($null, $null, 23.32),
    @($null, $null, 78.76),
    @($null, $null, 30.58),
    @($null, $null, 31.46),
    @($null, $null, 18.92)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $value = $row[$j]
        if ($null -ne $value) {
            $ws.Cells.Item($startRow + $i, $j + 1).Value = $value
        }
    }
}
